$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header for column F (copy formatting from E1, which shares the bold/bordered header style)
$ws.Range("F1").Value = "CAO_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Definition text for the first data row (BFO_0000029 / site)
$siteDef = "['B is a site means: b is a three-dimensional immaterial entity that is (partially or wholly) bounded by a material entity or it is a three-dimensional immaterial part thereof. [BFO]']"

$ws.Range("F2").Value = $siteDef
for ($r = 3; $r -le 15; $r++) {
    $ws.Range("F$r").Value = "[]"
}
